$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post at row 746 ("「コーヒーは繰り返さない恋物語」...") was removed.
# Deleting the entire row shifts every subsequent row up by one,
# matching the diff (row 747 -> 746, ..., row 807 -> 806) and updates
# the sheet's used-range dimension from A1:C807 to A1:C806 automatically.
$ws.Rows("746").Delete()
